$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-09-11 Wednesday" "2024-09-12 Thursday"

Replace-Text "681×7=4767" "207×8=1656"
Replace-Text "753×8=6024" "489×5=2445"
Replace-Text "320×4=1280" "369×8=2952"
Replace-Text "704×5=3520" "487×2=974"
Replace-Text "443×9=3987" "584×5=2920"

Replace-Text "938×5=4690" "177×8=1416"
Replace-Text "845×7=5915" "309×9=2781"
Replace-Text "998×2=1996" "607×6=3642"
Replace-Text "381×8=3048" "625×5=3125"
Replace-Text "624×6=3744" "636×3=1908"

Replace-Text "637×3=1911" "656×6=3936"
Replace-Text "128×7=896" "305×4=1220"
Replace-Text "551×6=3306" "266×7=1862"
Replace-Text "610×3=1830" "616×4=2464"
Replace-Text "845×5=4225" "219×2=438"

Replace-Text "379×6=2274" "785×8=6280"
Replace-Text "937×9=8433" "117×5=585"
Replace-Text "961×2=1922" "266×5=1330"
Replace-Text "901×5=4505" "838×8=6704"
Replace-Text "736×8=5888" "306×6=1836"

Replace-Text "846×3=2538" "968×5=4840"
Replace-Text "973×5=4865" "388×7=2716"
Replace-Text "997×7=6979" "824×2=1648"
Replace-Text "909×2=1818" "114×3=342"
Replace-Text "296×7=2072" "622×4=2488"

Write-Output "Done"
